# Scheduled data refresh: update leve-profit calculations (currentAveragePrice*,
# LevePrice*, LeveProfit*) per sheet/row from the latest market-board pull.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2750
$ws.Range("J40").Value = 2500
$ws.Range("L40").Value = 2500
$ws.Range("N40").Value = -2850
$ws.Range("H70").Value = 966.6667
$ws.Range("I70").Value = 633.3333
$ws.Range("K70").Value = 1899.9999
$ws.Range("M70").Value = -1629.9999
$ws.Range("H73").Value = 966.6667
$ws.Range("I73").Value = 633.3333
$ws.Range("K73").Value = 1899.9999
$ws.Range("M73").Value = -963.9999
$ws.Range("H132").Value = 1055.375
$ws.Range("I132").Value = 1057.1613
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 3171.4839
$ws.Range("L132").Value = 3000
$ws.Range("M132").Value = -641.4839000000002
$ws.Range("N132").Value = -8060
$ws.Range("H138").Value = 2272.53
$ws.Range("I138").Value = 1272.5217
$ws.Range("J138").Value = 2571.234
$ws.Range("K138").Value = 3817.5651
$ws.Range("L138").Value = 7713.701999999999
$ws.Range("M138").Value = 1322.4349
$ws.Range("N138").Value = -17993.702

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4460.85
$ws.Range("I32").Value = 4460.85
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 4460.85
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -4173.85
$ws.Range("N32").ClearContents()
$ws.Range("H102").Value = 1007.2308
$ws.Range("I102").Value = 1007.2308
$ws.Range("K102").Value = 1007.2308
$ws.Range("M102").Value = 614.7692
$ws.Range("H122").Value = 1593.1578
$ws.Range("I122").Value = 1551.2354
$ws.Range("J122").Value = 1949.5
$ws.Range("K122").Value = 4653.706200000001
$ws.Range("L122").Value = 5848.5
$ws.Range("M122").Value = -2203.706200000001
$ws.Range("N122").Value = -10748.5
$ws.Range("H132").Value = 1819.7241
$ws.Range("I132").Value = 1419.1613
$ws.Range("J132").Value = 2279.6296
$ws.Range("K132").Value = 4257.4839
$ws.Range("L132").Value = 6838.888800000001
$ws.Range("M132").Value = -1727.4839
$ws.Range("N132").Value = -11898.8888

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 6714.7144
$ws.Range("I12").Value = 5333.8335
$ws.Range("J12").Value = 15000
$ws.Range("K12").Value = 5333.8335
$ws.Range("L12").Value = 15000
$ws.Range("M12").Value = -5165.8335
$ws.Range("N12").Value = -15336
$ws.Range("H54").Value = 7950
$ws.Range("I54").Value = 5266.6665
$ws.Range("K54").Value = 5266.6665
$ws.Range("M54").Value = -4782.6665
$ws.Range("H94").Value = 292
$ws.Range("I94").Value = 322.66666
$ws.Range("K94").Value = 322.66666
$ws.Range("M94").Value = 128.33334
$ws.Range("H97").Value = 20100
$ws.Range("J97").Value = 22000
$ws.Range("L97").Value = 22000
$ws.Range("N97").Value = -23982
$ws.Range("H99").Value = 1870.6666
$ws.Range("I99").Value = 1804.8334
$ws.Range("K99").Value = 1804.8334
$ws.Range("M99").Value = -306.8334
$ws.Range("H134").Value = 6568.3335
$ws.Range("I134").Value = 7000
$ws.Range("J134").Value = 6352.5
$ws.Range("K134").Value = 21000
$ws.Range("L134").Value = 19057.5
$ws.Range("M134").Value = -18465
$ws.Range("N134").Value = -24127.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1399.5
$ws.Range("J16").Value = 1332
$ws.Range("L16").Value = 1332
$ws.Range("N16").Value = -1906
$ws.Range("H31").Value = 2582.5925
$ws.Range("I31").Value = 2418.625
$ws.Range("J31").Value = 2821.0908
$ws.Range("K31").Value = 2418.625
$ws.Range("L31").Value = 2821.0908
$ws.Range("M31").Value = -2123.625
$ws.Range("N31").Value = -3411.0908
$ws.Range("H34").Value = 2582.5925
$ws.Range("I34").Value = 2418.625
$ws.Range("J34").Value = 2821.0908
$ws.Range("K34").Value = 2418.625
$ws.Range("L34").Value = 2821.0908
$ws.Range("M34").Value = -2216.625
$ws.Range("N34").Value = -3225.0908
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()
$ws.Range("H99").Value = 627078.7
$ws.Range("I99").Value = 771243
$ws.Range("K99").Value = 771243
$ws.Range("M99").Value = -769745
$ws.Range("H113").Value = 1399.5
$ws.Range("J113").Value = 1332
$ws.Range("L113").Value = 1332
$ws.Range("N113").Value = -5672
$ws.Range("H126").Value = 627078.7
$ws.Range("I126").Value = 771243
$ws.Range("K126").Value = 2313729
$ws.Range("M126").Value = -2311259
$ws.Range("H134").Value = 1506.4546
$ws.Range("I134").Value = 1357.8889
$ws.Range("K134").Value = 4073.6667
$ws.Range("M134").Value = -1538.6667
$ws.Range("H141").Value = 71458.336
$ws.Range("J141").Value = 71350
$ws.Range("L141").Value = 71350
$ws.Range("N141").Value = -81710

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 114359.91
$ws.Range("I4").Value = 56776.668
$ws.Range("K4").Value = 170330.004
$ws.Range("M4").Value = -170218.004
$ws.Range("H6").Value = 767.5
$ws.Range("I6").Value = 933.3333
$ws.Range("K6").Value = 2799.9999
$ws.Range("M6").Value = -2686.9999
$ws.Range("H56").Value = 6570.913
$ws.Range("I56").Value = 6570.913
$ws.Range("K56").Value = 6570.913
$ws.Range("M56").Value = -6040.913
$ws.Range("H114").Value = 3197.2222
$ws.Range("I114").Value = 571.6667
$ws.Range("J114").Value = 4510
$ws.Range("K114").Value = 1715.0001
$ws.Range("L114").Value = 13530
$ws.Range("M114").Value = 1538.9999
$ws.Range("N114").Value = -20038
$ws.Range("H129").Value = 28215.074
$ws.Range("I129").Value = 765
$ws.Range("J129").Value = 30411.08
$ws.Range("K129").Value = 2295
$ws.Range("L129").Value = 91233.24000000001
$ws.Range("M129").Value = 2705
$ws.Range("N129").Value = -101233.24
$ws.Range("H131").Value = 6958016.5
$ws.Range("J131").Value = 14555.03
$ws.Range("L131").Value = 43665.09
$ws.Range("N131").Value = -53745.09
$ws.Range("H134").Value = 45747.13
$ws.Range("I134").Value = 54679.74
$ws.Range("J134").Value = 3317.25
$ws.Range("K134").Value = 164039.22
$ws.Range("L134").Value = 9951.75
$ws.Range("M134").Value = -158969.22
$ws.Range("N134").Value = -20091.75

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1559.8
$ws.Range("I122").Value = 1324.75
$ws.Range("K122").Value = 3974.25
$ws.Range("M122").Value = -1524.25
$ws.Range("H126").Value = 2780682.2
$ws.Range("I126").Value = 3270912.5
$ws.Range("J126").Value = 2711
$ws.Range("K126").Value = 9812737.5
$ws.Range("L126").Value = 8133
$ws.Range("M126").Value = -9810267.5
$ws.Range("N126").Value = -13073
$ws.Range("H132").Value = 1103969.1
$ws.Range("I132").Value = 1608412.8
$ws.Range("K132").Value = 4825238.4
$ws.Range("M132").Value = -4822708.4
$ws.Range("H134").Value = 25000
$ws.Range("J134").Value = 25000
$ws.Range("L134").Value = 75000
$ws.Range("N134").Value = -80070

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2951.5
$ws.Range("I7").Value = 2547.5557
$ws.Range("K7").Value = 2547.5557
$ws.Range("M7").Value = -2435.5557
$ws.Range("H40").Value = 9582.538
$ws.Range("J40").Value = 10438.8
$ws.Range("L40").Value = 10438.8
$ws.Range("N40").Value = -10710.8
$ws.Range("H126").Value = 2951.5
$ws.Range("I126").Value = 2547.5557
$ws.Range("K126").Value = 7642.6671
$ws.Range("M126").Value = -5172.6671
$ws.Range("H132").Value = 4138.846
$ws.Range("I132").Value = 2339.4
$ws.Range("K132").Value = 7018.200000000001
$ws.Range("M132").Value = -4488.200000000001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 53388.168
$ws.Range("J46").Value = 53388.168
$ws.Range("L46").Value = 53388.168
$ws.Range("N46").Value = -53850.168
$ws.Range("H100").Value = 459.63635
$ws.Range("I100").Value = 459.63635
$ws.Range("K100").Value = 919.2727
$ws.Range("M100").Value = -378.2727
$ws.Range("H108").Value = 79999
$ws.Range("J108").Value = 79999
$ws.Range("L108").Value = 79999
$ws.Range("N108").Value = -87679
$ws.Range("H126").Value = 5572.115
$ws.Range("I126").Value = 6231
$ws.Range("K126").Value = 18693
$ws.Range("M126").Value = -16223
$ws.Range("H132").Value = 1759.2745
$ws.Range("I132").Value = 1638.6666
$ws.Range("J132").Value = 2048.7334
$ws.Range("K132").Value = 4915.9998
$ws.Range("L132").Value = 6146.2002
$ws.Range("M132").Value = -2385.9998
$ws.Range("N132").Value = -11206.2002
$ws.Range("H134").Value = 53388.168
$ws.Range("J134").Value = 53388.168
$ws.Range("L134").Value = 160164.504
$ws.Range("N134").Value = -165234.504
